$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.988.73"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "3.847.59"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'703.67"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").Value = "'172.49"
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").Value = "3.845.84"
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("D11").Value = "'7.31"
$ws.Range("E11").Value = "  -2.40%  "
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").Value = "'36.32"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").Value = "4.490.80"
$ws.Range("D16").Value = "3.883.68"
$ws.Range("E16").Value = "  +1.81%  "
$ws.Range("D17").Value = "70.947.60"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").Value = "'17.47"
$ws.Range("E20").Value = "  -2.78%  "
$ws.Range("D21").Value = "'10.77"
$ws.Range("E21").Value = "  -3.97%  "
$ws.Range("D22").Value = "'493.70"
$ws.Range("E22").Value = "  +2.08%  "
$ws.Range("D23").Value = "'0.718"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'84.75"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("D25").Value = "'0.0000148"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").Value = "'12.16"
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("E28").Value = "  -3.13%  "
$ws.Range("D29").Value = "'3.16"
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("D33").Value = "'0.183"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("D34").Value = "'29.44"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").Value = "3.799.23"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").Value = "'9.15"
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").Value = "'2.38"
$ws.Range("E39").Value = "  +6.97%  "
$ws.Range("E40").Value = "  +6.89%  "
$ws.Range("D41").Value = "'6.02"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").Value = "'3.35"
$ws.Range("E42").Value = "  -5.34%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "'0.000315"
$ws.Range("E45").Value = "  -6.07%  "
$ws.Range("D46").Value = "'163.16"
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("D47").Value = "'48.85"
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("B48").Value = "TheGraph"
$ws.Range("C48").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D48").Value = "'0.299"
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'8.64"
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").Value = "'43.39"
$ws.Range("E50").Value = "  -3.93%  "
$ws.Range("D51").Value = "'408.55"
$ws.Range("E51").Value = "  +2.19%  "
